$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The "_GoBack" bookmark used to sit in the last cell of the big
#    results table at the bottom of the document; it is being moved
#    to the new credentials table inserted below, so strip it from
#    its old home first (keeping the now-empty paragraph/cell intact).
#    (Doing this before the structural table insert below avoids a
#    stale-Range duplication quirk when both edits touch tables.)
# ------------------------------------------------------------------
$bigTable = $d.Tables($d.Tables.Count)
$lastCell = $bigTable.Cell($bigTable.Rows.Count, $bigTable.Columns.Count)
$emptyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p>'
$lastCell.Range.InsertXML($emptyXml) | Out-Null

# ------------------------------------------------------------------
# 2) The "Resources" section's old placeholder paragraph ("	: -")
#    becomes a Heading2 bullet paragraph announcing the mock login
#    data, followed by a 2x2 table of credentials and a blank
#    paragraph.
# ------------------------------------------------------------------
$target = $d.Paragraphs(5).Range

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>&#8226; Mock data provided of admin account</w:t></w:r><w:r><w:t xml:space="preserve"> for login</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>' + `
          '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4675"/><w:gridCol w:w="4675"/></w:tblGrid>' + `
          '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>User-email</w:t></w:r></w:p></w:tc>' + `
          '<w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Password</w:t></w:r></w:p></w:tc></w:tr>' + `
          '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>admin1@cmuAD.ac.th</w:t></w:r></w:p></w:tc>' + `
          '<w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>1</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:tc></w:tr>' + `
          '</w:tbl>' + `
          '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'

$target.InsertXML($newXml) | Out-Null

Write-Output "ok"
